$d = $word.ActiveDocument

# The graduation date for Oregon State University changed from 6/2027 to 6/2026.
$d.Content.Find.Execute("6/2027", $true, $false, $false, $false, $false, $true, 1, $false, "6/2026", 2)
